$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume text cells keep their exact literal formatting (no numeric auto-conversion)
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '70.159.11'
$ws.Range('E2').Value = '  +1.67%  '
$ws.Range('D3').Value = '3.533.96'
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '594.98'
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('D6').Value = '170.92'
$ws.Range('E6').Value = '  +1.24%  '
$ws.Range('D7').Value = '3.531.07'
$ws.Range('E7').Value = '  +0.72%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('E10').Value = '  +3.49%  '
$ws.Range('D11').Value = '7.40'
$ws.Range('E11').Value = '  +9.64%  '
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('D13').Value = '46.12'
$ws.Range('E13').Value = '  -2.22%  '
$ws.Range('E14').Value = '  +0.35%  '
$ws.Range('D15').Value = '4.111.04'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = '8.28'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('D17').Value = '607.88'
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').Value = '3.534.81'
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('D19').Value = '70.241.31'
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').Value = '17.23'
$ws.Range('E21').Value = '  -0.81%  '
$ws.Range('D22').Value = '0.874'
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('D23').Value = '9.13'
$ws.Range('E23').Value = '  -17.28%  '
$ws.Range('D24').Value = '15.60'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').Value = '96.11'
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('D26').Value = '3.69'
$ws.Range('E26').Value = '  -3.38%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('D29').Value = '33.71'
$ws.Range('E29').Value = '  +3.62%  '
$ws.Range('D30').Value = '8.97'
$ws.Range('E30').Value = '  -1.71%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '715.50'
$ws.Range('E31').Value = '  +14.88%  '
$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').Value = '3.02'
$ws.Range('E32').Value = '  -2.82%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '8.17'
$ws.Range('E33').Value = '  -3.72%  '
$ws.Range('D34').Value = '7.04'
$ws.Range('E34').Value = '  +2.57%  '
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('D36').Value = '0.0999'
$ws.Range('E36').Value = '  -1.53%  '
$ws.Range('D37').Value = '3.54'
$ws.Range('E37').Value = '  +1.50%  '
$ws.Range('D38').Value = '10.70'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').Value = '0.0472'
$ws.Range('E39').Value = '  +7.38%  '
$ws.Range('D40').Value = '56.85'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').Value = '0.141'
$ws.Range('E42').Value = '  +4.61%  '
$ws.Range('D43').Value = '3.376.39'
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('E44').Value = '  -2.24%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').Value = '2.92'
$ws.Range('E45').Value = '  +6.81%  '
$ws.Range('B46').Value = 'PEPE'
$ws.Range('C46').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D46').Value = '0.0₃0693'
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '32.37'
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('D48').Value = '2.57'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('D50').Value = '132.63'
$ws.Range('E50').Value = '  -0.87%  '
$ws.Range('E51').Value = '  -0.10%  '
